$d = $word.ActiveDocument

# --- Append the new bulleted note at the end of the document -------------
# Continue the existing bulleted list (same numId as the paragraphs above
# it) by inserting a new paragraph after the last one and typing the text.
$end = $d.Content
$end.Collapse(0)             # wdCollapseEnd
$end.InsertParagraphAfter()
$end.Collapse(0)

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Har lavet lille firkant. Burde det være ens?"

# --- Register a second (unused) bullet-list definition --------------------
# While taking notes, a second bullet list gallery entry got created in the
# document (e.g. from re-applying "Bullets" while checking formatting) but
# the visible paragraph above kept using the original list (numId 1). Model
# that by applying the default bullet template to a throwaway paragraph,
# which mints a new abstractNum/num pair, then remove the throwaway text
# again without touching the real paragraph's numbering.
$scratchEnd = $d.Content
$scratchEnd.Collapse(0)
$scratchEnd.InsertParagraphAfter()
$scratchEnd.Collapse(0)
$scratchPara = $d.Paragraphs.Last
$scratchPara.Range.Text = "x"

$bulletTemplate = $word.ListGalleries.Item(1).ListTemplates.Item(1)
$scratchPara.Range.ListFormat.ApplyListTemplate($bulletTemplate)

$scratchPara.Range.Delete()
